# WRN8IMB remaining code
# Update the single construction-site data row (row 2 of the "ConstructionSite"
# table) with the new 8IMB test site reference for 2024-09-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Test8IMB20240918SITE63629"
